# "Add score over time tracker" -- append newly solved problems (all
# solved on 2026-01-25) to the bottom of the "Solution Dates" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newProblems = @(
    "administrativeproblems",
    "babynames",
    "bst",
    "caching",
    "continuousmedian",
    "cookieselection",
    "doctorkattis",
    "kattissquest",
    "orphanbackups"
)
$solvedDate = "2026-01-25"

$startRow = 194
for ($i = 0; $i -lt $newProblems.Count; $i++) {
    $row = $startRow + $i

    # Column B holds dates formatted as plain text (e.g. "2026-01-24" on
    # the existing rows) rather than real date serials, so force a text
    # number format before assigning the value -- otherwise Excel's
    # auto-detection would silently convert the string into a date.
    $ws.Range("B$row").NumberFormat = "@"

    $ws.Range("A$row").Value = $newProblems[$i]
    $ws.Range("B$row").Value = $solvedDate
}
